# Scheduled market-data refresh: recompute currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) columns (H:N) for the leves whose
# underlying item prices moved since the last run. Columns A:G are untouched.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 26: Everything Is Impossible / Budding Ash Wand
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

# Row 32: Automata for the People / Crab Oil
$ws.Range("H32").Value = 3833.6667
$ws.Range("I32").Value = 5333.3335
$ws.Range("J32").Value = 2334
$ws.Range("K32").Value = 5333.3335
$ws.Range("L32").Value = 2334
$ws.Range("M32").Value = -5007.3335
$ws.Range("N32").Value = -2986

# Row 111: An Eye for Healing / Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 3217711.8
$ws.Range("I111").Value = 2196.7
$ws.Range("J111").Value = 8576903
$ws.Range("K111").Value = 6590.099999999999
$ws.Range("L111").Value = 25730709
$ws.Range("M111").Value = -3523.099999999999
$ws.Range("N111").Value = -25736843

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 2591.5
$ws.Range("I112").Value = 3283.3333
$ws.Range("J112").Value = 2360.889
$ws.Range("K112").Value = 9849.999899999999
$ws.Range("L112").Value = 7082.667
$ws.Range("M112").Value = -8741.999899999999
$ws.Range("N112").Value = -9298.667000000001

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 10000
$ws.Range("N113").Value = -16508
$ws.Range("M113").ClearContents()

# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 13081.963
$ws.Range("I125").Value = 27768.584
$ws.Range("J125").Value = 1332.6666
$ws.Range("K125").Value = 249917.256
$ws.Range("L125").Value = 11993.9994
$ws.Range("M125").Value = -247457.256
$ws.Range("N125").Value = -16913.9994

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 63886.375
$ws.Range("I137").Value = 1547.4
$ws.Range("J137").Value = 108414.21
$ws.Range("K137").Value = 4642.200000000001
$ws.Range("L137").Value = 325242.63
$ws.Range("M137").Value = -2092.200000000001
$ws.Range("N137").Value = -330342.63

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 5116.8887
$ws.Range("I138").Value = 2451
$ws.Range("J138").Value = 7249.6
$ws.Range("K138").Value = 7353
$ws.Range("L138").Value = 21748.8
$ws.Range("M138").Value = -2213
$ws.Range("N138").Value = -32028.8

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24: A Firm Hand / Iron Gauntlets
$ws.Range("H24").Value = 29016.25
$ws.Range("J24").Value = 29016.25
$ws.Range("L24").Value = 29016.25
$ws.Range("N24").Value = -29764.25

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 127053.57
$ws.Range("I32").Value = 135244.95
$ws.Range("J32").Value = 93543.37
$ws.Range("K32").Value = 135244.95
$ws.Range("L32").Value = 93543.37
$ws.Range("M32").Value = -134957.95
$ws.Range("N32").Value = -94117.37

# Row 55: Employee Retention / Mythril Elmo
$ws.Range("H55").Value = 20000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

# Row 80: A Squire to Inspire / Titanium Hoplon
$ws.Range("H80").Value = 19997.143
$ws.Range("J80").Value = 19997.143
$ws.Range("L80").Value = 19997.143
$ws.Range("N80").Value = -21993.143

# Row 83: All's Fair in Highborn Assassination (L) / Titanium Hoplon
$ws.Range("H83").Value = 19997.143
$ws.Range("J83").Value = 19997.143
$ws.Range("L83").Value = 59991.429
$ws.Range("N83").Value = -69975.429

# Row 96: The Gauntlet Is Cast / High Steel Gauntlets of Fending
$ws.Range("H96").Value = 60085.75
$ws.Range("J96").Value = 60085.75
$ws.Range("L96").Value = 60085.75
$ws.Range("N96").Value = -65577.75

# Row 100: En Garde and on Guard / Doman Iron Gauntlets of Fending
$ws.Range("H100").Value = 29016.25
$ws.Range("J100").Value = 29016.25
$ws.Range("L100").Value = 29016.25
$ws.Range("N100").Value = -31180.25

# Row 112: Wrapped Knuckles / Deepgold Gloves of Fending
$ws.Range("H112").Value = 70033.5
$ws.Range("J112").Value = 70033.5
$ws.Range("L112").Value = 70033.5
$ws.Range("N112").Value = -72987.5

# Row 114: A New Regular / Bluespirit Gauntlets of Fending
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 13702.842
$ws.Range("I122").Value = 16383.667
$ws.Range("J122").Value = 3649.75
$ws.Range("K122").Value = 49151.001
$ws.Range("L122").Value = 10949.25
$ws.Range("M122").Value = -46701.001
$ws.Range("N122").Value = -15849.25

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 17337.63
$ws.Range("I132").Value = 18258.883
$ws.Range("J132").Value = 9507
$ws.Range("K132").Value = 54776.649
$ws.Range("L132").Value = 28521
$ws.Range("M132").Value = -52246.649
$ws.Range("N132").Value = -33581

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2758.7942
$ws.Range("I31").Value = 1868.1892
$ws.Range("J31").Value = 3821.7742
$ws.Range("K31").Value = 1868.1892
$ws.Range("L31").Value = 3821.7742
$ws.Range("M31").Value = -1573.1892
$ws.Range("N31").Value = -4411.7742

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2758.7942
$ws.Range("I34").Value = 1868.1892
$ws.Range("J34").Value = 3821.7742
$ws.Range("K34").Value = 1868.1892
$ws.Range("L34").Value = 3821.7742
$ws.Range("M34").Value = -1666.1892
$ws.Range("N34").Value = -4225.7742

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 4451.316
$ws.Range("I132").Value = 4143.0557
$ws.Range("K132").Value = 12429.1671
$ws.Range("M132").Value = -9899.167099999999

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 1928386.1
$ws.Range("J68").Value = 2783529.5
$ws.Range("L68").Value = 8350588.5
$ws.Range("N68").Value = -8352210.5

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 1928386.1
$ws.Range("J71").Value = 2783529.5
$ws.Range("L71").Value = 25051765.5
$ws.Range("N71").Value = -25059877.5

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 128546.266
$ws.Range("I131").Value = 56703.89
$ws.Range("J131").Value = 290191.62
$ws.Range("K131").Value = 170111.67
$ws.Range("L131").Value = 870574.86
$ws.Range("M131").Value = -165071.67
$ws.Range("N131").Value = -880654.86

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 20: Brothers in Arms / Brass Wristlets of Crafting
$ws.Range("H20").Value = 7800
$ws.Range("I20").Value = 7500
$ws.Range("J20").Value = 8000
$ws.Range("K20").Value = 7500
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = -7255
$ws.Range("N20").Value = -8490

# Row 24: Bad Guys Eat Brass / Brass Ring of Crafting
$ws.Range("H24").Value = 1503635.8
$ws.Range("I24").Value = 10501450
$ws.Range("J24").Value = 4000
$ws.Range("K24").Value = 10501450
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = -10501277
$ws.Range("N24").Value = -4346

# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 22274
$ws.Range("I97").Value = 34884.367
$ws.Range("K97").Value = 34884.367
$ws.Range("M97").Value = -34388.367

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 2074.291
$ws.Range("I102").Value = 1896.5532
$ws.Range("J102").Value = 3118.5
$ws.Range("K102").Value = 1896.5532
$ws.Range("L102").Value = 3118.5
$ws.Range("M102").Value = -274.5532000000001
$ws.Range("N102").Value = -6362.5

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1989.0476
$ws.Range("I122").Value = 2109.7576
$ws.Range("J122").Value = 1546.4445
$ws.Range("K122").Value = 6329.2728
$ws.Range("L122").Value = 4639.333500000001
$ws.Range("M122").Value = -3879.2728
$ws.Range("N122").Value = -9539.333500000001

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 41710736
$ws.Range("I61").Value = 47621790
$ws.Range("J61").Value = 333333
$ws.Range("K61").Value = 47621790
$ws.Range("L61").Value = 333333
$ws.Range("M61").Value = -47621588
$ws.Range("N61").Value = -333737

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 41710736
$ws.Range("I113").Value = 47621790
$ws.Range("J113").Value = 333333
$ws.Range("K113").Value = 47621790
$ws.Range("L113").Value = 333333
$ws.Range("M113").Value = -47619620
$ws.Range("N113").Value = -337673

# Row 114: A Heady Endeavor / Atrociraptorskin Headgear of Scouting
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3790.5625
$ws.Range("I122").Value = 3588.4614
$ws.Range("J122").Value = 4666.3335
$ws.Range("K122").Value = 10765.3842
$ws.Range("L122").Value = 13999.0005
$ws.Range("M122").Value = -8315.3842
$ws.Range("N122").Value = -18899.0005

# Row 134: Freezing Fingers / Crocodileskin Fingerless Gloves of Striking
$ws.Range("H134").Value = 55214.5
$ws.Range("J134").Value = 55214.5
$ws.Range("L134").Value = 55214.5
$ws.Range("N134").Value = -65354.5

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 10040.167
$ws.Range("I113").Value = 3949.25
$ws.Range("K113").Value = 11847.75
$ws.Range("M113").Value = -9677.75
